$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "68.595.04"
# or "0.0000290" keep their exact original formatting instead of being
# auto-parsed into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '68.595.04'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.866.55'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '603.14'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').Value = '173.50'
$ws.Range('E6').Value = '  +4.30%  '
$ws.Range('D7').Value = '3.869.21'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('E10').Value = '  +3.08%  '
$ws.Range('D11').Value = '6.55'
$ws.Range('E11').Value = '  +3.69%  '
$ws.Range('D12').Value = '0.0000290'
$ws.Range('E12').Value = '  +16.63%  '
$ws.Range('D13').Value = '0.463'
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range('D14').Value = '37.39'
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').Value = '4.515.59'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').Value = '3.854.22'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '68.658.94'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').Value = '18.39'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').Value = '7.50'
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').Value = '10.93'
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('D22').Value = '473.12'
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('D23').Value = '0.737'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').Value = '84.07'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('D27').Value = '12.30'
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('D28').Value = '10.56'
$ws.Range('E28').Value = '  +5.80%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('D31').Value = '4.018.24'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').Value = '7.81'
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').Value = '31.33'
$ws.Range('E34').Value = '  +1.06%  '
$ws.Range('D35').Value = '9.49'
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').Value = '3.833.42'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '3.99'
$ws.Range('E37').Value = '  +21.96%  '
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('D39').Value = '6.02'
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  +3.58%  '
$ws.Range('E44').Value = '  +11.20%  '
$ws.Range('D45').Value = '2.00'
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('B46').Value = 'Cosmos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D46').Value = '8.84'
$ws.Range('E46').Value = '  +3.68%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '423.66'
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('D49').Value = '46.71'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('D50').Value = '0.0362'
$ws.Range('E50').Value = '  +2.70%  '
$ws.Range('D51').Value = '142.41'
$ws.Range('E51').Value = '  -0.77%  '

# Remove the temporary text number-format so the cells end up with no
# explicit style, matching the style-less cells in the original sheet.
$ws.Range("D2:D51").ClearFormats()
